# Updated OCR app features and fixes
# Adds new ultrasound/patient-detail columns (AL:AW) and four new report rows (8-11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header columns AL1:AW1
# ---------------------------------------------------------------------------
$newHeaders = @(
    "Patient Name",
    "Patient Age",
    "Patient Gender",
    "Liver Size",
    "Gall Bladder Status",
    "Spleen Size",
    "Pancreas Status",
    "Right Kidney Size",
    "Left Kidney Size",
    "Urinary Bladder Status",
    "Ultrasound Findings",
    "Ultrasound Impression"
)

# AL = column 38
$startCol = 38
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $newHeaders[$i]
}

# Match the bold/centered/bordered header style used by the existing header
# row (A1:AK1) instead of leaving the new header cells unstyled.
$ws.Range("A1").Copy()
$ws.Range("AL1:AW1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Blank out the new columns (AL:AW) for the pre-existing rows 2-7 so the
#    used range keeps every column touched like the rest of the sheet.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    for ($c = 38; $c -le 49; $c++) {
        $ws.Cells.Item($r, $c).Value = ""
    }
}

# The text columns below can contain values that look like dates/numbers
# ("2025-12-27", "90", ...). Force those specific cells to Text format
# *before* assigning so Excel's COM layer doesn't silently reinterpret them
# as date serials / numbers (matching how the source file stored them).
$ws.Range("A8:A11").NumberFormat = "@"
$ws.Cells.Item(11, 39).NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3. Row 8: duplicate of row 7's lab values, new date + report type
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "2025-12-27"
$ws.Cells.Item(8, 2).Value = "Liver Function Test (LFT)"
$ws.Cells.Item(8, 3).Value = 11.3
$ws.Cells.Item(8, 4).Value = 4.45
$ws.Cells.Item(8, 5).Value = 10.54
$ws.Cells.Item(8, 6).Value = 416
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = ""
$ws.Cells.Item(8, 9).Value = 13
$ws.Cells.Item(8, 10).Value = 9
$ws.Cells.Item(8, 11).Value = ""
$ws.Cells.Item(8, 12).Value = ""
$ws.Cells.Item(8, 13).Value = ""
$ws.Cells.Item(8, 14).Value = 0.4
$ws.Cells.Item(8, 15).Value = 0.2
$ws.Cells.Item(8, 16).Value = 0.2
$ws.Cells.Item(8, 17).Value = 34
$ws.Cells.Item(8, 18).Value = 5
$ws.Cells.Item(8, 19).Value = 219
$ws.Cells.Item(8, 20).Value = ""
$ws.Cells.Item(8, 21).Value = 4.6
$ws.Cells.Item(8, 22).Value = 1.9
$ws.Cells.Item(8, 23).Value = 2.4
$ws.Cells.Item(8, 24).Value = 34.7
$ws.Cells.Item(8, 25).Value = 78
$ws.Cells.Item(8, 26).Value = 25.4
$ws.Cells.Item(8, 27).Value = 32
$ws.Cells.Item(8, 28).Value = 12
$ws.Cells.Item(8, 29).Value = 8.199999999999999
$ws.Cells.Item(8, 30).Value = ""
$ws.Cells.Item(8, 31).Value = ""
$ws.Cells.Item(8, 32).Value = ""
$ws.Cells.Item(8, 33).Value = 26
$ws.Cells.Item(8, 34).Value = 79
$ws.Cells.Item(8, 35).Value = 2
$ws.Cells.Item(8, 36).Value = 2
$ws.Cells.Item(8, 37).Value = 10
for ($c = 38; $c -le 49; $c++) {
    $ws.Cells.Item(8, $c).Value = ""
}

# ---------------------------------------------------------------------------
# 4. Row 9: Vitals Check (hip-ultrasound style findings)
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = "2025-12-28"
$ws.Cells.Item(9, 2).Value = "Vitals Check"
for ($c = 3; $c -le 37; $c++) {
    $ws.Cells.Item(9, $c).Value = ""
}
$ws.Cells.Item(9, 38).Value = ""
$ws.Cells.Item(9, 39).Value = 1
$ws.Cells.Item(9, 40).Value = "Gender Leger"
$ws.Cells.Item(9, 41).Value = 68
$ws.Cells.Item(9, 42).Value = "Normal"
$ws.Cells.Item(9, 43).Value = 48
$ws.Cells.Item(9, 44).Value = ""
$ws.Cells.Item(9, 45).Value = "43 x 19 mm"
$ws.Cells.Item(9, 46).Value = "48 x 22 mm"
$ws.Cells.Item(9, 47).Value = ""
$avText = @"
Acetabulum and iliac wing appear normal.
Femoral head appears normal and cartilagenous.
Central ossific nucleus not visualised in femoral head on either side.
No joint effusion.
No intramusclar fluid collection.
Right angle Alpha - 61
Beta - 56
Left angle Alpha - 62
Beta - 55
No obvious e/o developmental dysplasia of hip.
"@
$avText = $avText.TrimEnd("`r", "`n")
$awText = @"
;
Normal study. .
"@
$awText = $awText.TrimEnd("`r", "`n")
$ws.Cells.Item(9, 48).Value = $avText
$ws.Cells.Item(9, 49).Value = $awText

# ---------------------------------------------------------------------------
# 5. Row 10: Ultrasound Report (same patient-detail values as row 9)
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = "2025-12-28"
$ws.Cells.Item(10, 2).Value = "Ultrasound Report"
for ($c = 3; $c -le 37; $c++) {
    $ws.Cells.Item(10, $c).Value = ""
}
$ws.Cells.Item(10, 38).Value = ""
$ws.Cells.Item(10, 39).Value = 1
$ws.Cells.Item(10, 40).Value = "Gender Leger"
$ws.Cells.Item(10, 41).Value = 68
$ws.Cells.Item(10, 42).Value = "Normal"
$ws.Cells.Item(10, 43).Value = 48
$ws.Cells.Item(10, 44).Value = ""
$ws.Cells.Item(10, 45).Value = "43 x 19 mm"
$ws.Cells.Item(10, 46).Value = "48 x 22 mm"
$ws.Cells.Item(10, 47).Value = ""
$ws.Cells.Item(10, 48).Value = $avText
$ws.Cells.Item(10, 49).Value = $awText

# ---------------------------------------------------------------------------
# 6. Row 11: Thyroid Test
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "2025-12-28"
$ws.Cells.Item(11, 2).Value = "Thyroid Test"
for ($c = 3; $c -le 29; $c++) {
    $ws.Cells.Item(11, $c).Value = ""
}
$ws.Cells.Item(11, 30).Value = 176.2
$ws.Cells.Item(11, 31).Value = 10.9
$ws.Cells.Item(11, 32).Value = 2.98
for ($c = 33; $c -le 38; $c++) {
    $ws.Cells.Item(11, $c).Value = ""
}
$ws.Cells.Item(11, 39).Value = "90"
for ($c = 40; $c -le 49; $c++) {
    $ws.Cells.Item(11, $c).Value = ""
}
